# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): rows keyed by F-cell address -> new value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 3412
$wsExhibit.Range("F5").Value  = 6971
$wsExhibit.Range("F6").Value  = 2439
$wsExhibit.Range("F7").Value  = 42
$wsExhibit.Range("F8").Value  = 110
$wsExhibit.Range("F12").Value = 33
$wsExhibit.Range("F13").Value = 173

# Sheet "全部类型" (sheet4.xml): same events, shifted by one extra row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 3412
$wsAll.Range("F6").Value  = 6971
$wsAll.Range("F7").Value  = 2439
$wsAll.Range("F8").Value  = 42
$wsAll.Range("F9").Value  = 110
$wsAll.Range("F13").Value = 33
$wsAll.Range("F14").Value = 173
